$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H9").Value2 = 322.625
$ws.Range("I9").Value2 = 85.75
$ws.Range("J9").Value2 = 559.5
$ws.Range("K9").Value2 = 85.75
$ws.Range("L9").Value2 = 559.5
$ws.Range("M9").Value2 = 83.25
$ws.Range("N9").Value2 = -897.5
$ws.Range("H62").Value2 = 9342.714
$ws.Range("J62").Value2 = 4999
$ws.Range("L62").Value2 = 4999
$ws.Range("N62").Value2 = -6247
$ws.Range("H64").Value2 = 4437.5
$ws.Range("I64").Value2 = 4333.3335
$ws.Range("J64").Value2 = 4500
$ws.Range("K64").Value2 = 4333.3335
$ws.Range("L64").Value2 = 4500
$ws.Range("M64").Value2 = -4085.3335
$ws.Range("N64").Value2 = -4996
$ws.Range("H65").Value2 = 9342.714
$ws.Range("J65").Value2 = 4999
$ws.Range("L65").Value2 = 24995
$ws.Range("N65").Value2 = -31235
$ws.Range("H67").Value2 = 4437.5
$ws.Range("I67").Value2 = 4333.3335
$ws.Range("J67").Value2 = 4500
$ws.Range("K67").Value2 = 4333.3335
$ws.Range("L67").Value2 = 4500
$ws.Range("M67").Value2 = -3475.3335
$ws.Range("N67").Value2 = -6216
$ws.Range("H74").Value2 = 2399
$ws.Range("I74").Value2 = 2399
$ws.Range("K74").Value2 = 2399
$ws.Range("M74").Value2 = -1463
$ws.Range("H77").Value2 = 2399
$ws.Range("I77").Value2 = 2399
$ws.Range("K77").Value2 = 11995
$ws.Range("M77").Value2 = -7315
$ws.Range("H106").Value2 = 985
$ws.Range("I106").Value2 = 985
$ws.Range("K106").Value2 = 985
$ws.Range("M106").Value2 = -354
$ws.Range("H112").Value2 = 2390.8333
$ws.Range("J112").Value2 = 3054.923
$ws.Range("L112").Value2 = 9164.769
$ws.Range("N112").Value2 = -11380.769
$ws.Range("H115").Value2 = 1006.93335
$ws.Range("I115").Value2 = 828.8570999999999
$ws.Range("K115").Value2 = 2486.5713
$ws.Range("M115").Value2 = -919.5712999999996
$ws.Range("H127").Value2 = 2887.1304
$ws.Range("I127").Value2 = 3121.85
$ws.Range("J127").Value2 = 1322.3334
$ws.Range("K127").Value2 = 9365.549999999999
$ws.Range("L127").Value2 = 3967.0002
$ws.Range("M127").Value2 = -4405.549999999999
$ws.Range("N127").Value2 = -13887.0002
$ws.Range("H129").Value2 = 499247
$ws.Range("I129").Value2 = 581871.5
$ws.Range("J129").Value2 = 3500
$ws.Range("K129").Value2 = 1745614.5
$ws.Range("L129").Value2 = 10500
$ws.Range("M129").Value2 = -1740614.5
$ws.Range("N129").Value2 = -20500
$ws.Range("H131").Value2 = 1373.2
$ws.Range("I131").Value2 = 815.2308
$ws.Range("K131").Value2 = 2445.6924
$ws.Range("M131").Value2 = 2594.3076
$ws.Range("H133").Value2 = 69265.45
$ws.Range("J133").Value2 = 69265.45
$ws.Range("L133").Value2 = 69265.45
$ws.Range("N133").Value2 = -79385.45
$ws.Range("H141").Value2 = 7217.5
$ws.Range("J141").Value2 = 9537.77
$ws.Range("L141").Value2 = 28613.31
$ws.Range("N141").Value2 = -38973.31
$ws = $wb.Worksheets.Item(2)
$ws.Range("H5").Value2 = 7575
$ws.Range("J5").Value2 = 15000
$ws.Range("L5").Value2 = 15000
$ws.Range("N5").Value2 = -15224
$ws.Range("H31").Value2 = 3690
$ws.Range("I31").Value2 = 3690
$ws.Range("J31").Value2 = 0
$ws.Range("K31").Value2 = 3690
$ws.Range("L31").Value2 = 0
$ws.Range("M31").Value2 = -3396
$ws.Range("N31").Value2 = $null
$ws.Range("H32").Value2 = 145506.53
$ws.Range("I32").Value2 = 157645.61
$ws.Range("K32").Value2 = 157645.61
$ws.Range("M32").Value2 = -157358.61
$ws.Range("H41").Value2 = 2468
$ws.Range("I41").Value2 = 1908.8889
$ws.Range("J41").Value2 = 7500
$ws.Range("K41").Value2 = 1908.8889
$ws.Range("L41").Value2 = 7500
$ws.Range("M41").Value2 = -1494.8889
$ws.Range("N41").Value2 = -8328
$ws.Range("H45").Value2 = 4999.6665
$ws.Range("I45").Value2 = 4999.6665
$ws.Range("K45").Value2 = 4999.6665
$ws.Range("M45").Value2 = -4622.6665
$ws.Range("H61").Value2 = 5699.5386
$ws.Range("I61").Value2 = 7066
$ws.Range("J61").Value2 = 3513.2
$ws.Range("K61").Value2 = 7066
$ws.Range("L61").Value2 = 3513.2
$ws.Range("M61").Value2 = -6854
$ws.Range("N61").Value2 = -3937.2
$ws.Range("H97").Value2 = 2205.6191
$ws.Range("I97").Value2 = 1184.3889
$ws.Range("K97").Value2 = 1184.3889
$ws.Range("M97").Value2 = -688.3888999999999
$ws.Range("H110").Value2 = 1098.2916
$ws.Range("I110").Value2 = 877.0526
$ws.Range("K110").Value2 = 877.0526
$ws.Range("M110").Value2 = 1167.9474
$ws.Range("H122").Value2 = 2113.5
$ws.Range("J122").Value2 = 1570.6666
$ws.Range("L122").Value2 = 4711.9998
$ws.Range("N122").Value2 = -9611.9998
$ws.Range("H132").Value2 = 25000000
$ws.Range("I132").Value2 = 25000000
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 75000000
$ws.Range("L132").Value2 = 0
$ws.Range("M132").Value2 = -74997470
$ws.Range("N132").Value2 = $null
$ws.Range("H133").Value2 = 76326.664
$ws.Range("J133").Value2 = 76326.664
$ws.Range("L133").Value2 = 76326.664
$ws.Range("N133").Value2 = -81386.664
$ws.Range("H134").Value2 = 65578.8
$ws.Range("J134").Value2 = 65578.8
$ws.Range("L134").Value2 = 65578.8
$ws.Range("N134").Value2 = -75718.8
$ws.Range("H136").Value2 = 5699.5386
$ws.Range("I136").Value2 = 7066
$ws.Range("J136").Value2 = 3513.2
$ws.Range("K136").Value2 = 21198
$ws.Range("L136").Value2 = 10539.6
$ws.Range("M136").Value2 = -18648
$ws.Range("N136").Value2 = -15639.6
$ws = $wb.Worksheets.Item(3)
$ws.Range("H4").Value2 = 7575
$ws.Range("J4").Value2 = 15000
$ws.Range("L4").Value2 = 15000
$ws.Range("N4").Value2 = -15230
$ws.Range("H12").Value2 = 306.5
$ws.Range("I12").Value2 = 410
$ws.Range("J12").Value2 = 203
$ws.Range("K12").Value2 = 410
$ws.Range("L12").Value2 = 203
$ws.Range("M12").Value2 = -242
$ws.Range("N12").Value2 = -539
$ws.Range("H17").Value2 = 463.6
$ws.Range("I17").Value2 = 100
$ws.Range("J17").Value2 = 554.5
$ws.Range("K17").Value2 = 100
$ws.Range("L17").Value2 = 554.5
$ws.Range("M17").Value2 = 72
$ws.Range("N17").Value2 = -898.5
$ws.Range("H25").Value2 = 1362.4
$ws.Range("I25").Value2 = 1253
$ws.Range("J25").Value2 = 1800
$ws.Range("K25").Value2 = 1253
$ws.Range("L25").Value2 = 1800
$ws.Range("M25").Value2 = -1018
$ws.Range("N25").Value2 = -2270
$ws.Range("H86").Value2 = 1869.1666
$ws.Range("I86").Value2 = 1665.6666
$ws.Range("J86").Value2 = 2276.1667
$ws.Range("K86").Value2 = 1665.6666
$ws.Range("L86").Value2 = 2276.1667
$ws.Range("M86").Value2 = -542.6666
$ws.Range("N86").Value2 = -4522.1667
$ws.Range("H89").Value2 = 1869.1666
$ws.Range("I89").Value2 = 1665.6666
$ws.Range("J89").Value2 = 2276.1667
$ws.Range("K89").Value2 = 8328.333000000001
$ws.Range("L89").Value2 = 11380.8335
$ws.Range("M89").Value2 = -2712.333000000001
$ws.Range("N89").Value2 = -22612.8335
$ws.Range("H94").Value2 = 2624.9062
$ws.Range("I94").Value2 = 2286.1155
$ws.Range("J94").Value2 = 4093
$ws.Range("K94").Value2 = 2286.1155
$ws.Range("L94").Value2 = 4093
$ws.Range("M94").Value2 = -1835.1155
$ws.Range("N94").Value2 = -4995
$ws.Range("H105").Value2 = 2675.875
$ws.Range("I105").Value2 = 2036.5294
$ws.Range("K105").Value2 = 2036.5294
$ws.Range("M105").Value2 = -289.5293999999999
$ws.Range("H107").Value2 = 1091
$ws.Range("I107").Value2 = 1049.5
$ws.Range("J107").Value2 = 1153.25
$ws.Range("K107").Value2 = 1049.5
$ws.Range("L107").Value2 = 1153.25
$ws.Range("M107").Value2 = 870.5
$ws.Range("N107").Value2 = -4993.25
$ws.Range("H134").Value2 = 11233.5
$ws.Range("I134").Value2 = 4355.75
$ws.Range("K134").Value2 = 13067.25
$ws.Range("M134").Value2 = -10532.25
$ws = $wb.Worksheets.Item(4)
$ws.Range("H15").Value2 = 0
$ws.Range("I15").Value2 = 0
$ws.Range("K15").Value2 = 0
$ws.Range("M15").Value2 = $null
$ws.Range("H22").Value2 = 2082.1875
$ws.Range("I22").Value2 = 903.4286
$ws.Range("K22").Value2 = 903.4286
$ws.Range("M22").Value2 = -553.4286
$ws.Range("H56").Value2 = 50093
$ws.Range("I56").Value2 = 50093
$ws.Range("K56").Value2 = 50093
$ws.Range("M56").Value2 = -49248
$ws.Range("H99").Value2 = 14205500
$ws.Range("I99").Value2 = 86999.60000000001
$ws.Range("K99").Value2 = 86999.60000000001
$ws.Range("M99").Value2 = -85501.60000000001
$ws.Range("H126").Value2 = 14205500
$ws.Range("I126").Value2 = 86999.60000000001
$ws.Range("K126").Value2 = 260998.8
$ws.Range("M126").Value2 = -258528.8
$ws.Range("H132").Value2 = 2345.15
$ws.Range("I132").Value2 = 1898.5
$ws.Range("J132").Value2 = 4131.75
$ws.Range("K132").Value2 = 5695.5
$ws.Range("L132").Value2 = 12395.25
$ws.Range("M132").Value2 = -3165.5
$ws.Range("N132").Value2 = -17455.25
$ws.Range("H134").Value2 = 6669568.5
$ws.Range("I134").Value2 = 8573431
$ws.Range("K134").Value2 = 25720293
$ws.Range("M134").Value2 = -25717758
$ws = $wb.Worksheets.Item(5)
$ws.Range("H9").Value2 = 281.83334
$ws.Range("I9").Value2 = 122.75
$ws.Range("J9").Value2 = 600
$ws.Range("K9").Value2 = 368.25
$ws.Range("L9").Value2 = 1800
$ws.Range("M9").Value2 = -144.25
$ws.Range("N9").Value2 = -2248
$ws.Range("H22").Value2 = 337.5
$ws.Range("I22").Value2 = 175.5
$ws.Range("J22").Value2 = 499.5
$ws.Range("K22").Value2 = 526.5
$ws.Range("L22").Value2 = 1498.5
$ws.Range("M22").Value2 = -357.5
$ws.Range("N22").Value2 = -1836.5
$ws.Range("H23").Value2 = 99.333336
$ws.Range("I23").Value2 = 51.5
$ws.Range("J23").Value2 = 123.25
$ws.Range("K23").Value2 = 154.5
$ws.Range("L23").Value2 = 369.75
$ws.Range("M23").Value2 = 80.5
$ws.Range("N23").Value2 = -839.75
$ws.Range("H27").Value2 = 337.5
$ws.Range("I27").Value2 = 175.5
$ws.Range("J27").Value2 = 499.5
$ws.Range("K27").Value2 = 526.5
$ws.Range("L27").Value2 = 1498.5
$ws.Range("M27").Value2 = -424.5
$ws.Range("N27").Value2 = -1702.5
$ws.Range("H31").Value2 = 0
$ws.Range("J31").Value2 = 0
$ws.Range("L31").Value2 = 0
$ws.Range("N31").Value2 = $null
$ws.Range("H58").Value2 = 6563.636
$ws.Range("I58").Value2 = 1066.6666
$ws.Range("J58").Value2 = 8625
$ws.Range("K58").Value2 = 3199.9998
$ws.Range("L58").Value2 = 25875
$ws.Range("M58").Value2 = -3071.9998
$ws.Range("N58").Value2 = -26131
$ws.Range("H109").Value2 = 4784.5293
$ws.Range("J109").Value2 = 9375
$ws.Range("L109").Value2 = 28125
$ws.Range("N109").Value2 = -30205
$ws.Range("H113").Value2 = 1474.1428
$ws.Range("J113").Value2 = 1604.6666
$ws.Range("L113").Value2 = 4813.9998
$ws.Range("N113").Value2 = -9153.9998
$ws.Range("H131").Value2 = 7555.1113
$ws.Range("J131").Value2 = 7555.1113
$ws.Range("L131").Value2 = 22665.3339
$ws.Range("N131").Value2 = -32745.3339
$ws.Range("H134").Value2 = 6231.591
$ws.Range("J134").Value2 = 14000
$ws.Range("L134").Value2 = 42000
$ws.Range("N134").Value2 = -52140
$ws.Range("H140").Value2 = 1207.8
$ws.Range("I140").Value2 = 1207.8
$ws.Range("K140").Value2 = 3623.4
$ws.Range("M140").Value2 = 1556.6
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value2 = 76.73333
$ws.Range("I2").Value2 = 46.666668
$ws.Range("K2").Value2 = 46.666668
$ws.Range("M2").Value2 = 66.333332
$ws.Range("H14").Value2 = 143718060
$ws.Range("I14").Value2 = 143718060
$ws.Range("J14").Value2 = 0
$ws.Range("K14").Value2 = 143718060
$ws.Range("L14").Value2 = 0
$ws.Range("M14").Value2 = -143717892
$ws.Range("N14").Value2 = $null
$ws.Range("H17").Value2 = 226
$ws.Range("J17").Value2 = 207.5
$ws.Range("L17").Value2 = 207.5
$ws.Range("N17").Value2 = -543.5
$ws.Range("H19").Value2 = 33500
$ws.Range("I19").Value2 = 30000
$ws.Range("K19").Value2 = 30000
$ws.Range("M19").Value2 = -29712
$ws.Range("H21").Value2 = 82000
$ws.Range("J21").Value2 = 151250
$ws.Range("L21").Value2 = 151250
$ws.Range("N21").Value2 = -151596
$ws.Range("H22").Value2 = 1142.5
$ws.Range("I22").Value2 = 285
$ws.Range("J22").Value2 = 2000
$ws.Range("K22").Value2 = 285
$ws.Range("L22").Value2 = 2000
$ws.Range("M22").Value2 = 244
$ws.Range("N22").Value2 = -3058
$ws.Range("H23").Value2 = 6734.6665
$ws.Range("J23").Value2 = 7575
$ws.Range("L23").Value2 = 7575
$ws.Range("N23").Value2 = -8021
$ws.Range("H25").Value2 = 2399.8
$ws.Range("J25").Value2 = 2499.75
$ws.Range("L25").Value2 = 2499.75
$ws.Range("N25").Value2 = -3557.75
$ws.Range("H30").Value2 = 82000
$ws.Range("J30").Value2 = 151250
$ws.Range("L30").Value2 = 151250
$ws.Range("N30").Value2 = -151460
$ws.Range("H122").Value2 = 2906.32
$ws.Range("I122").Value2 = 2271.5264
$ws.Range("K122").Value2 = 6814.5792
$ws.Range("M122").Value2 = -4364.5792
$ws.Range("H126").Value2 = 7732.0713
$ws.Range("I126").Value2 = 8437.416999999999
$ws.Range("K126").Value2 = 25312.251
$ws.Range("M126").Value2 = -22842.251
$ws.Range("H132").Value2 = 16029.75
$ws.Range("I132").Value2 = 16913.566
$ws.Range("J132").Value2 = 2772.5
$ws.Range("K132").Value2 = 50740.698
$ws.Range("L132").Value2 = 8317.5
$ws.Range("M132").Value2 = -48210.698
$ws.Range("N132").Value2 = -13377.5
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value2 = 2006.909
$ws.Range("I16").Value2 = 1841.5555
$ws.Range("J16").Value2 = 2751
$ws.Range("K16").Value2 = 1841.5555
$ws.Range("L16").Value2 = 2751
$ws.Range("M16").Value2 = -1671.5555
$ws.Range("N16").Value2 = -3091
$ws.Range("H22").Value2 = 3190.5293
$ws.Range("J22").Value2 = 3692.1428
$ws.Range("L22").Value2 = 3692.1428
$ws.Range("N22").Value2 = -4282.1428
$ws.Range("H27").Value2 = 3190.5293
$ws.Range("J27").Value2 = 3692.1428
$ws.Range("L27").Value2 = 3692.1428
$ws.Range("N27").Value2 = -3906.1428
$ws.Range("H46").Value2 = 3772.0344
$ws.Range("H61").Value2 = 19885.182
$ws.Range("I61").Value2 = 16637.666
$ws.Range("J61").Value2 = 34499
$ws.Range("K61").Value2 = 16637.666
$ws.Range("L61").Value2 = 34499
$ws.Range("M61").Value2 = -16435.666
$ws.Range("N61").Value2 = -34903
$ws.Range("H111").Value2 = 0
$ws.Range("J111").Value2 = 0
$ws.Range("L111").Value2 = 0
$ws.Range("N111").Value2 = $null
$ws.Range("H113").Value2 = 19885.182
$ws.Range("I113").Value2 = 16637.666
$ws.Range("J113").Value2 = 34499
$ws.Range("K113").Value2 = 16637.666
$ws.Range("L113").Value2 = 34499
$ws.Range("M113").Value2 = -14467.666
$ws.Range("N113").Value2 = -38839
$ws.Range("H132").Value2 = 3574321.8
$ws.Range("I132").Value2 = 7144914
$ws.Range("K132").Value2 = 21434742
$ws.Range("M132").Value2 = -21432212
$ws = $wb.Worksheets.Item(8)
$ws.Range("H29").Value2 = 8704.083000000001
$ws.Range("J29").Value2 = 7500
$ws.Range("L29").Value2 = 7500
$ws.Range("N29").Value2 = -8080
$ws.Range("H100").Value2 = 1397.4
$ws.Range("I100").Value2 = 1441.6111
$ws.Range("K100").Value2 = 2883.2222
$ws.Range("M100").Value2 = -2342.2222
$ws.Range("H113").Value2 = 465
$ws.Range("I113").Value2 = 449.46155
$ws.Range("K113").Value2 = 1348.38465
$ws.Range("M113").Value2 = 821.61535
$ws.Range("H122").Value2 = 32375.055
$ws.Range("I122").Value2 = 1124.6072
$ws.Range("K122").Value2 = 3373.8216
$ws.Range("M122").Value2 = -923.8215999999998
$ws.Range("H126").Value2 = 1843
$ws.Range("I126").Value2 = 2076.1875
$ws.Range("K126").Value2 = 6228.5625
$ws.Range("M126").Value2 = -3758.5625
$ws.Range("H132").Value2 = 1726.7693
$ws.Range("I132").Value2 = 1737.579
$ws.Range("J132").Value2 = 1697.4286
$ws.Range("K132").Value2 = 5212.737
$ws.Range("L132").Value2 = 5092.2858
$ws.Range("M132").Value2 = -2682.737
$ws.Range("N132").Value2 = -10152.2858
$ws.Range("H136").Value2 = 305729.66
$ws.Range("I136").Value2 = 355852.1
$ws.Range("K136").Value2 = 1067556.3
$ws.Range("M136").Value2 = -1065006.3
